# Update countries & provincias Spain
# Refresh the "last updated" timestamp and the per-country COVID figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 17:19"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5159325
$ws.Range("C4").Value = 9602
$ws.Range("D4").Value = 2638800
$ws.Range("E4").Value = 2355339
$ws.Range("G4").Value = 116
$ws.Range("H4").Value = 165186

# --- Row 6: India ---
$ws.Range("B6").Value = 2199101
$ws.Range("C6").Value = 47081
$ws.Range("D6").Value = 1506413
$ws.Range("E6").Value = 648640
$ws.Range("G6").Value = 595
$ws.Range("H6").Value = 44048

# --- Row 15: Reino Unido ---
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 46574

# --- Row 19: Italia ---
$ws.Range("B19").Value = 250566
$ws.Range("C19").Value = 463
$ws.Range("D19").Value = 202098
$ws.Range("E19").Value = 13263
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 35205

# --- Row 48: Portugal ---
$ws.Range("B48").Value = 52668
$ws.Range("C48").Value = 131
$ws.Range("D48").Value = 38511
$ws.Range("E48").Value = 12401
$ws.Range("G48").Value = 6
$ws.Range("H48").Value = 1756

# --- Row 64: Moldavia ---
$ws.Range("B64").Value = 27660
$ws.Range("C64").Value = 217
$ws.Range("D64").Value = 19300
$ws.Range("E64").Value = 7515
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 845

# --- Row 102: Grecia ---
$ws.Range("D102").Value = 3804
$ws.Range("E102").Value = 1406

# --- Row 133: Islandia ---
$ws.Range("B133").Value = 1958
$ws.Range("C133").Value = 3
$ws.Range("D133").Value = 1834
$ws.Range("E133").Value = 114

# --- Row 178: Trinidad yTobago ---
$ws.Range("B178").Value = 279
$ws.Range("C178").Value = 4
$ws.Range("E178").Value = 136
